$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.345812797546387
$ws.Range("B1").Value = 4.082015037536621
$ws.Range("C1").Value = 2.069024801254272
$ws.Range("D1").Value = 1.582843065261841
$ws.Range("E1").Value = 1.421786546707153
